$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '87.924.67'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.75%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.249.18'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -3.22%  '

# Row 4
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '212.45'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -4.45%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '626.85'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -3.52%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.385'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +12.05%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.711'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +15.01%  '

# Row 9
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.04%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '3.246.15'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -3.42%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.575'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -6.65%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.189'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +12.69%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000268'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -3.08%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.49'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +0.07%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '34.15'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -3.53%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.852.66'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -3.31%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '87.860.70'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.52%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.272.14'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -2.34%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.23'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.72%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.04'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -5.00%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '436.10'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -7.12%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '8.97'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -3.62%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.33'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -5.03%  '

# Row 24
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.24%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.35'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -2.84%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '12.47'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -7.76%  '

# Row 27
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +10.78%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '3.420.03'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -2.50%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '77.27'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -2.70%  '

# Row 30
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.01%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.176'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -11.98%  '

# Row 32
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.74%  '

# Row 33
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '569.48'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -5.08%  '

# Row 34
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '8.86'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -6.01%  '

# Row 35
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.39'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -11.06%  '

# Row 36
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '7.27'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +4.50%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.96'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -5.53%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.139'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -8.61%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '22.93'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -5.35%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.28'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +6.19%  '

# Row 41
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.09%  '

# Row 42
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '21.81'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +0.55%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.403'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -5.06%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.03'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -6.39%  '

# Row 45
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.03%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '151.57'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -3.98%  '

# Row 47
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '180.00'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -7.12%  '

# Row 48
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.135'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +17.95%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '45.22'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -5.28%  '

# Row 50
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -3.43%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.25'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -2.64%  '
